# Major Updates 4th April 2016
# Add a new course/teacher entry (row 61) to the timetable sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A61").Value = "CS666"
$ws.Range("B61").Value = "Test Course"
$ws.Range("C61").Value = "CS3"
$ws.Range("D61").Value = "test111"

# Move the selection to the newly entered cell, matching the author's
# final cursor position when the workbook was saved.
[void]$ws.Range("C61").Select()
